$d = $word.ActiveDocument

# Locate the two target paragraphs by content rather than a hard-coded index,
# so the script is resilient to minor paragraph-numbering differences.
$para1 = $null
$para2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text
    if ($para1 -eq $null -and $ptext -like "*32bit (see 4), the do the following:*") {
        $para1 = $d.Paragraphs($i)
    }
    if ($para2 -eq $null -and $ptext -like "*setup_env32*" -and $ptext -like "*press enter*") {
        $para2 = $d.Paragraphs($i)
    }
}

# --- Fix 1: typo "...32bit (see 4), the do the following:" -> "...32bit (see 4), then do the following:"
#     Insert the missing "n" that turns "the" into "then" - a minimal, surgical
#     edit that leaves the rest of the paragraph (and the _GoBack bookmark)
#     untouched.
if ($para1 -ne $null) {
    $fr1 = $para1.Range.Duplicate
    $found1 = $fr1.Find.Execute(", the do the following:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    Write-Host "Fix1 found: $found1"
    if ($found1) {
        $insertPoint = $d.Range($fr1.Start + 5, $fr1.Start + 5)
        $insertPoint.InsertBefore("n")
    }
} else {
    Write-Host "Fix1: target paragraph not found"
}

# --- Fix 2: color the closing-quote tail ", and press enter" blue (0070C0) on
#     the "Type "python setup_env32.py gsflow", and press enter" line.
if ($para2 -ne $null) {
    $rng2 = $para2.Range.Duplicate
    $found2 = $rng2.Find.Execute([char]0x201D + ", and press enter", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    Write-Host "Fix2 found: $found2 text=[$($rng2.Text)]"
    if ($found2) {
        $rng2.Font.Color = 12611584
    }
} else {
    Write-Host "Fix2: target paragraph not found"
}
